$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 312 (Terminal La Palmera de
# La Serena - Papa). Insert a whole new row there so every existing record
# from the old row 312 onward shifts down by one (old row 377 -> new 378),
# then populate the freshly inserted row with its data.
$ws.Rows.Item(312).Insert()

$ws.Cells.Item(312, 1).Value = 8
$ws.Cells.Item(312, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(312, 3).Value = "Coquimbo"
$ws.Cells.Item(312, 4).Value = 44637
$ws.Cells.Item(312, 5).Value = 4
$ws.Cells.Item(312, 6).Value = 100114001
$ws.Cells.Item(312, 7).Value = "Papa"
$ws.Cells.Item(312, 8).Value = "Asterix"
$ws.Cells.Item(312, 9).Value = "1a (cosecha)"
$ws.Cells.Item(312, 10).Value = 2000
$ws.Cells.Item(312, 11).Value = 8000
$ws.Cells.Item(312, 12).Value = 9000
$ws.Cells.Item(312, 13).Value = 8500
$ws.Cells.Item(312, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(312, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(312, 16).Value = 340
$ws.Cells.Item(312, 17).Value = 25
$ws.Cells.Item(312, 18).Value = "Hortaliza"
